$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Force the cell to store a literal text value even when the string
    # looks like a number (e.g. "600.48"), matching the workbook's
    # original inlineStr/text encoding for the Price column, then restore
    # the default "Normal" style so no stray number-format style lingers.
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "74.766.89"
$ws.Range("E2").Value = "  +1.12%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.841.13"
$ws.Range("E3").Value = "  +10.43%  "

# Row 4 - TetherUSD (price unchanged)
$ws.Range("E4").Value = "  -0.02%  "

# Row 5 - BNB
Set-TextValue "D5" "600.48"
$ws.Range("E5").Value = "  +4.23%  "

# Row 6 - Solana
Set-TextValue "D6" "188.12"
$ws.Range("E6").Value = "  +2.22%  "

# Row 7 - USDC (price unchanged)
$ws.Range("E7").Value = "  -0.11%  "

# Row 8 - XRP
Set-TextValue "D8" "0.555"
$ws.Range("E8").Value = "  +4.25%  "

# Row 9 - Dogecoin
Set-TextValue "D9" "0.193"
$ws.Range("E9").Value = "  -3.82%  "

# Row 10 - LidoStakedEther
Set-TextValue "D10" "2.841.57"
$ws.Range("E10").Value = "  +10.51%  "

# Row 11 - TRON (price unchanged)
$ws.Range("E11").Value = "  -0.21%  "

# Row 12 - Cardano
Set-TextValue "D12" "0.371"
$ws.Range("E12").Value = "  +3.72%  "

# Row 13 - Toncoin
Set-TextValue "D13" "4.88"
$ws.Range("E13").Value = "  +2.69%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextValue "D14" "3.363.89"
$ws.Range("E14").Value = "  +10.34%  "

# Row 15 - WrappedBTC
Set-TextValue "D15" "74.791.70"
$ws.Range("E15").Value = "  +1.06%  "

# Row 16 & 17 swap: Avalanche now ranked 16, ShibaInu now ranked 17
$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue "D16" "27.17"
$ws.Range("E16").Value = "  +4.28%  "

$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D17" "0.0000188"
$ws.Range("E17").Value = "  -0.86%  "

# Row 18 - WrappedEther
Set-TextValue "D18" "2.842.06"
$ws.Range("E18").Value = "  +10.19%  "

# Row 19 - Uniswap
Set-TextValue "D19" "9.14"
$ws.Range("E19").Value = "  +7.27%  "

# Row 20 - Chainlink
Set-TextValue "D20" "12.49"
$ws.Range("E20").Value = "  +7.19%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "374.79"
$ws.Range("E21").Value = "  -0.35%  "

# Row 22 - SuiNetwork
Set-TextValue "D22" "2.26"
$ws.Range("E22").Value = "  -1.57%  "

# Row 23 - Polkadot
Set-TextValue "D23" "4.13"
$ws.Range("E23").Value = "  +2.14%  "

# Row 24 - LEO
Set-TextValue "D24" "6.18"
$ws.Range("E24").Value = "  -0.47%  "

# Row 25 - Dai (price unchanged)
$ws.Range("E25").Value = "  -0.01%  "

# Row 26 - Litecoin
Set-TextValue "D26" "70.83"
$ws.Range("E26").Value = "  +2.21%  "

# Row 27 - NEARProtocol
Set-TextValue "D27" "4.21"
$ws.Range("E27").Value = "  +1.85%  "

# Row 28 - WrappedeETH (price unchanged)
$ws.Range("E28").Value = "  +9.83%  "

# Row 29 - Aptos
Set-TextValue "D29" "9.54"
$ws.Range("E29").Value = "  +4.17%  "

# Row 30 - PEPE (price unchanged)
$ws.Range("E30").Value = "  +11.52%  "

# Row 31 - Binance-PegBSC-USD (price unchanged)
$ws.Range("E31").Value = "  -0.14%  "

# Row 32 - Bittensor
Set-TextValue "D32" "524.41"
$ws.Range("E32").Value = "  +5.57%  "

# Row 33 - Fetch.AI
Set-TextValue "D33" "1.40"
$ws.Range("E33").Value = "  +5.40%  "

# Row 34 - InternetComputer(DFINITY)
Set-TextValue "D34" "7.93"
$ws.Range("E34").Value = "  +0.60%  "

# Row 35 - PancakeSwap (price unchanged)
$ws.Range("E35").Value = "  +6.40%  "

# Row 36 - FirstDigitalUSD (price unchanged)
$ws.Range("E36").Value = "  -0.04%  "

# Row 37 - Kaspa (price unchanged)
$ws.Range("E37").Value = "  +2.25%  "

# Row 38 - EthereumClassic
Set-TextValue "D38" "20.13"
$ws.Range("E38").Value = "  +5.29%  "

# Row 39 - Monero
Set-TextValue "D39" "162.58"
$ws.Range("E39").Value = "  +2.32%  "

# Row 40 - WhiteBITCoin
Set-TextValue "D40" "19.28"
$ws.Range("E40").Value = "  -0.31%  "

# Row 41 - Aave
Set-TextValue "D41" "185.27"
$ws.Range("E41").Value = "  +25.13%  "

# Row 42 - USDe (price unchanged)
$ws.Range("E42").Value = "  +0.00%  "

# Row 43 - RenderToken
Set-TextValue "D43" "5.07"
$ws.Range("E43").Value = "  +4.05%  "

# Row 44 - PolygonEcosystemToken (price unchanged)
$ws.Range("E44").Value = "  +6.86%  "

# Row 45 - Stacks (price unchanged)
$ws.Range("E45").Value = "  +2.12%  "

# Row 46 - ImmutableX
Set-TextValue "D46" "1.24"
$ws.Range("E46").Value = "  +8.41%  "

# Row 47 & 48 swap: OKB now ranked 47, dogwifhat now ranked 48
$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D47" "39.70"
$ws.Range("E47").Value = "  +2.00%  "

$ws.Range("B48").Value = "dogwifhat"
$ws.Range("C48").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D48" "2.36"
$ws.Range("E48").Value = "  -1.91%  "

# Row 49 - Cronos
Set-TextValue "D49" "0.0856"
$ws.Range("E49").Value = "  +4.74%  "

# Row 50 - ARBITRUM
Set-TextValue "D50" "0.573"
$ws.Range("E50").Value = "  +10.93%  "

# Row 51 - Filecoin (price unchanged)
$ws.Range("E51").Value = "  +4.51%  "
